$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-25, columns B,C,D,E,F,I,J,K,L,M,N (column index order below)
$data = @(
@(1.02,1.025287330836978,1.030465776721949,1.028914531964831,1.023784189986237,1.032515387450662,1.030457194955088,1.033276708156985,1.031729963737781,1.026614611126441,1.014176909803406),
@(1.02,1.026166700969666,1.031124049102863,1.029741124877531,1.025298023884957,1.032706435531005,1.030976196468509,1.033743744424053,1.032364542180733,1.027933472828708,1.014350681409861),
@(1.02,1.02673586711347,1.031550076379495,1.030276490438439,1.026277740133164,1.032828930143585,1.031311555537793,1.034045346811891,1.03277501599953,1.028786542648162,1.014462915354156),
@(1.02,1.026975181211041,1.031729196771104,1.030501678210586,1.02668965650833,1.03288015699531,1.031452427762252,1.034171996043313,1.032947544997897,1.029145099742343,1.014510048570951),
@(1.02,1.02701536523793,1.031759272933817,1.030539495249293,1.026758821765275,1.03288874237571,1.031476074209365,1.034193252524772,1.032976511356928,1.029205298853292,1.014517959512643),
@(1.02,1.026739064696673,1.031552469723253,1.030279498937521,1.026283244003014,1.032829615700056,1.031313438322679,1.034047039673142,1.03277732147541,1.028791333989754,1.014463545346913),
@(1.02,1.025584485769755,1.030688225653635,1.029193777976344,1.024295765761419,1.03258018609424,1.03063269100284,1.033434669479947,1.031944451378765,1.027060395432746,1.014235679683374),
@(1.02,1.023551174572408,1.029165975964639,1.027284495901501,1.020794643215804,1.032132044351775,1.02942954348563,1.03235100956728,1.030475769379538,1.024007633940238,1.013832565914037),
@(1.02,1.022196465849242,1.028151635842156,1.026014302260145,1.018461010386046,1.031827505009199,1.028625051330736,1.0316255120709,1.029495962100299,1.021970466819431,1.013562766303514),
@(1.02,1.02161006331086,1.0277125424091,1.025464933599919,1.017450567795399,1.031694268279657,1.028276132318078,1.031310643323756,1.029071536876947,1.021087827411172,1.01344569103889),
@(1.02,1.021392276864082,1.027549462730734,1.025260969388868,1.017075244552749,1.031644572691573,1.028146442986949,1.031193578680814,1.028913862403502,1.020759891362581,1.013402166587142),
@(1.02,1.021438991454686,1.027584442995497,1.025304716061406,1.017155752741248,1.03165524186205,1.028174265645158,1.031218694342798,1.028947685196397,1.020830238652454,1.013411504424315),
@(1.02,1.021592060410797,1.027699061796488,1.025448071897572,1.017419543486978,1.031690164618699,1.028265413896214,1.031300968932666,1.029058503944512,1.021060721864131,1.013442094058186),
@(1.02,1.02168637517105,1.027769684786072,1.025536410867392,1.01758207353074,1.031711654451452,1.028321562046778,1.031351646672632,1.029126779855667,1.02120271885079,1.013460936377302),
@(1.02,1.022235387588181,1.028180779674188,1.026050775447308,1.018528070439396,1.031836318659522,1.028648195969727,1.031646393656971,1.029524126414622,1.022029033088729,1.013570530940305),
@(1.02,1.022579821015603,1.028438681983647,1.026373592796934,1.019121475193425,1.031914150774395,1.028852932677753,1.031831087176648,1.029773328214729,1.022547212446304,1.01363920979404),
@(1.02,1.022780742032927,1.028589123809222,1.02656194778668,1.019467601956717,1.03195941681235,1.028972297263979,1.031938745948273,1.029918667736748,1.022849406609804,1.013679244854283),
@(1.02,1.02284925409391,1.028640422557345,1.02662618231497,1.019585623044197,1.031974828946851,1.029012988191365,1.031975442952046,1.029968222089782,1.022952438479962,1.013692891674994),
@(1.02,1.022542864619474,1.028411010283404,1.026338951227168,1.019057808122576,1.031905813789368,1.028830972041052,1.031811278530653,1.029746592838605,1.022491622021287,1.013631843700447),
@(1.02,1.021546984600399,1.027665308881012,1.025405854528491,1.017341863726,1.031679886409684,1.028238575360874,1.031276744105491,1.029025871246556,1.020992852670643,1.013433087207335),
@(1.02,1.020921006473985,1.027196568117061,1.024819733980232,1.016262979490509,1.031536647824457,1.027865618758127,1.03094003425971,1.028572585998354,1.020050025945034,1.013307904285077),
@(1.02,1.02125283302276,1.027445045661065,1.025130394830335,1.016834918375663,1.03161269396252,1.028063376819084,1.031118589826243,1.02881289415213,1.020549884258946,1.013374286636405),
@(1.02,1.022559563565861,1.028423513897285,1.02635460407264,1.019086576518365,1.03190958132224,1.028840895276969,1.031820229420868,1.029758673452641,1.022516741098273,1.013635172195688),
@(1.02,1.024076689220818,1.029559430746506,1.027777625031273,1.021699668433681,1.032248919200967,1.029741009020861,1.032631702391773,1.030855581650605,1.024797180551787,1.013936967445471)
)

$colIndexes = @(2,3,4,5,6,9,10,11,12,13,14)

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowValues = $data[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $colIndexes.Count; $c++) {
        $colIndex = $colIndexes[$c]
        $ws.Cells.Item($excelRow, $colIndex).Value = $rowValues[$c]
    }
}
